$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Existing RW06 rows 107:117 - the "RegressionTest" (D) column flips from
#    Yes -> No (matching the pattern used by every other completed form
#    block), and picks up the same cell style already used one row above
#    (D106, style used for "No" cells) instead of the stray style it had.
# ---------------------------------------------------------------------------
$ws.Range("D106").Copy() | Out-Null
$ws.Range("D107:D117").PasteSpecial(-4122) | Out-Null
$ws.Range("D107:D117").Value = "No"

# ---------------------------------------------------------------------------
# 2) New "probateFormsRWxx.feature" test cases -> rows 118:123
# ---------------------------------------------------------------------------

# Row 118 - same shape as the rows above it (single-line C cell)
$ws.Range("A107:E107").Copy() | Out-Null
$ws.Range("A118:E118").PasteSpecial(-4122) | Out-Null

# Row 119 - long wrapped C cell, taller row (matches rows 51/54/57/61)
$ws.Range("A51:E51").Copy() | Out-Null
$ws.Range("A119:E119").PasteSpecial(-4122) | Out-Null
$ws.Rows("119").RowHeight = 28

# Rows 120:123 - same shape as row 118
$ws.Range("A107:E107").Copy() | Out-Null
$ws.Range("A120:E123").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

$ws.Range("A118").Value = "TC_117"
$ws.Range("B118").Value = "probateFormsRWxx.feature"
$ws.Range("C118").Value = "Open Estate"
$ws.Range("D118").Value = "Yes"
$ws.Range("E118").Value = "No"

$ws.Range("A119").Value = "TC_118"
$ws.Range("B119").Value = "probateFormsRWxx.feature"
$ws.Range("C119").Value = "Verify that the county, estate name, and ""Also Known As"" (AKA) values are auto-populated from the selected estate."
$ws.Range("D119").Value = "Yes"
$ws.Range("E119").Value = "No"

$ws.Range("A120").Value = "TC_119"
$ws.Range("B120").Value = "probateFormsRWxx.feature"
$ws.Range("C120").Value = "Verify, text can be entered in all the text areas."
$ws.Range("D120").Value = "Yes"
$ws.Range("E120").Value = "No"

$ws.Range("A121").Value = "TC_120"
$ws.Range("B121").Value = "probateFormsRWxx.feature"
$ws.Range("C121").Value = "Verify, the name entered in 1st text area is reflected in the signature."
$ws.Range("D121").Value = "Yes"
$ws.Range("E121").Value = "No"

$ws.Range("A122").Value = "TC_121"
$ws.Range("B122").Value = "probateFormsRWxx.feature"
$ws.Range("C122").Value = "Verify that changes in the witness name field are reflected under the signature line and vice-versa."
$ws.Range("D122").Value = "Yes"
$ws.Range("E122").Value = "No"

$ws.Range("A123").Value = "TC_122"
$ws.Range("B123").Value = "probateFormsRWxx.feature"
$ws.Range("C123").Value = "Reset the RWxx form"
$ws.Range("D123").Value = "Yes"
$ws.Range("E123").Value = "No"

# ---------------------------------------------------------------------------
# 3) Data validation - extend the "Yes,No" list validation to the new D:E
#    cells on rows 118:123 (mirrors what Excel does automatically when you
#    fill a validated column down into freshly used rows).
# ---------------------------------------------------------------------------
$newRng = $ws.Range("D118:E123")
$newRng.Validation.Delete() | Out-Null
$newRng.Validation.Add(3, 1, 1, '"Yes,No"') | Out-Null

# ---------------------------------------------------------------------------
# 4) View state - scroll / selection ends up near the new last row, like a
#    user who just finished typing the new test cases.
# ---------------------------------------------------------------------------
$ws.Application.GoTo($ws.Range("C121"), $true)
$ws.Range("C121").Select() | Out-Null
